# Tambah referensi manual: 123123
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A79").Value = 78
$ws.Range("B79").Value = "'015"
$ws.Range("C79").Value = "Kementerian Keuangan"
$ws.Range("D79").Value = "'123123"
$ws.Range("E79").Value = "TEST"
$ws.Range("F79").Value = "TESTTTTTTTT AJAAAAAAAAAA"
